$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as text (preserve numeric-looking strings exactly)
function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# --- Price (column D) updates for unchanged-label rows ---
Set-TextCell $ws "D2" "248.98"
Set-TextCell $ws "D3" "21.81"
Set-TextCell $ws "D4" "5.344"
Set-TextCell $ws "D5" "0.05621"
Set-TextCell $ws "D6" "3.405"
Set-TextCell $ws "D7" "6.392"
Set-TextCell $ws "D8" "0.8177"
Set-TextCell $ws "D9" "0.9581"
Set-TextCell $ws "D10" "0.1413"
Set-TextCell $ws "D11" "0.07521"
Set-TextCell $ws "D12" "0.03176"
Set-TextCell $ws "D13" "0.03050"
Set-TextCell $ws "D14" "0.09343"
Set-TextCell $ws "D15" "3.567"
Set-TextCell $ws "D16" "0.001605"
Set-TextCell $ws "D17" "0.04701"
Set-TextCell $ws "D25" "0.3257"
Set-TextCell $ws "D28" "0.0003097"
Set-TextCell $ws "D40" "0.03955"
Set-TextCell $ws "D41" "0.007072"
Set-TextCell $ws "D42" "0.1062"
Set-TextCell $ws "D43" "0.003101"
Set-TextCell $ws "D44" "0.008682"
Set-TextCell $ws "D45" "0.00005811"
Set-TextCell $ws "D46" "0.00000000749"
Set-TextCell $ws "D47" "0.0005491"
Set-TextCell $ws "D48" "0.6793"
Set-TextCell $ws "D49" "0.1601"
Set-TextCell $ws "D50" "0.00002098"
Set-TextCell $ws "D51" "0.01009"

# --- Rows 18-24: coin list shifted by one position (row relabeling) ---
Set-TextCell $ws "B18" "TigerCash"
Set-TextCell $ws "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D18" "0.006410"
Set-TextCell $ws "E18" "17TigerCashTCH"
Set-TextCell $ws "B19" "HotbitToken"
Set-TextCell $ws "C19" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell $ws "D19" "0.005077"
Set-TextCell $ws "E19" "18HotbitTokenHTB"
Set-TextCell $ws "B20" "BitKan"
Set-TextCell $ws "C20" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell $ws "D20" "0.001034"
Set-TextCell $ws "E20" "19BitKanKAN"
Set-TextCell $ws "B21" "NitroEx"
Set-TextCell $ws "C21" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell $ws "D21" "0.0001498"
Set-TextCell $ws "E21" "20NitroExNTX"
Set-TextCell $ws "B22" "LEO"
Set-TextCell $ws "C22" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D22" "3.746"
Set-TextCell $ws "E22" "21LEOLEO"
Set-TextCell $ws "B23" "BTSEToken"
Set-TextCell $ws "C23" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell $ws "D23" "2.122"
Set-TextCell $ws "E23" "22BTSETokenBTSE"
Set-TextCell $ws "B24" "One"
Set-TextCell $ws "C24" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws "D24" "0.01149"
Set-TextCell $ws "E24" "23OneONEBestin24h"

# --- Misc label-only update ---
Set-TextCell $ws "E49" "48BOLOBOLO"
